# Trade #68 closed at 2026-02-18 00:24:32 - unknown UNKNOWN +0.000%
#
# This script applies the bookkeeping updates that follow from closing an
# existing open EMAArbitrage trade (early exit, ~breakeven) and opening a
# brand-new MarketMaking trade, mirrored across the "All Trades" roll-up
# sheet and the per-strategy "MarketMaking" / "EMAArbitrage" sheets, plus
# the aggregate figures on "Summary" and "Strategy Status".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - aggregate counters
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 96      # Total Trades: 95 -> 96
$wsSummary.Range("B9").Value = 46.88   # Win Rate %: 47.37 -> 46.88

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - EMAArbitrage row (row 2)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D2").Value = 7        # Trades: 6 -> 7
$wsStatus.Range("G2").Value = 42.86    # Win Rate %: 50 -> 42.86

# ---------------------------------------------------------------------
# 3) All Trades sheet - close out trade in row 97 (the EMAArbitrage trade
#    that now reads CLOSED / early_exit) and append the newly opened
#    MarketMaking trade as row 126.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Range("G97").Value = 0.01
$wsAll.Range("H97").Value = "CLOSED"
$wsAll.Range("K97").Value = 100.27
$wsAll.Range("L97").Value = "early_exit"
$wsAll.Range("M97").Value = 0.09

# New row 126 - force text columns to literal text so date-/number-like
# strings ("2026-02-18") aren't coerced into date serials.
$wsAll.Range("B126").NumberFormat = "@"
$wsAll.Range("C126").NumberFormat = "@"

$wsAll.Range("A126").Value = 125
$wsAll.Range("B126").Value = "2026-02-18"
$wsAll.Range("C126").Value = "00:24:26"
$wsAll.Range("D126").Value = "MarketMaking"
$wsAll.Range("E126").Value = "UP"
$wsAll.Range("F126").Value = 0.01
$wsAll.Range("H126").Value = "OPEN"
$wsAll.Range("I126").Value = 0
$wsAll.Range("J126").Value = 0
$wsAll.Range("K126").Value = 99.40967800952272
$wsAll.Range("M126").Value = 0
$wsAll.Range("N126").Value = 0
$wsAll.Range("O126").Value = 0
$wsAll.Range("P126").Value = 0.6
$wsAll.Range("Q126").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# 4) MarketMaking sheet - append the same newly opened trade as row 46
#    (column layout differs from "All Trades": Entry Slippage, Exit
#    Slippage, Confidence, Entry Reason, Exit Reason, Duration).
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

$wsMM.Range("B46").NumberFormat = "@"
$wsMM.Range("C46").NumberFormat = "@"

$wsMM.Range("A46").Value = 125
$wsMM.Range("B46").Value = "2026-02-18"
$wsMM.Range("C46").Value = "00:24:26"
$wsMM.Range("D46").Value = "MarketMaking"
$wsMM.Range("E46").Value = "UP"
$wsMM.Range("F46").Value = 0.01
$wsMM.Range("H46").Value = "OPEN"
$wsMM.Range("I46").Value = 0
$wsMM.Range("J46").Value = 0
$wsMM.Range("K46").Value = 99.40967800952272
$wsMM.Range("L46").Value = 0
$wsMM.Range("M46").Value = 0
$wsMM.Range("N46").Value = 0.6
$wsMM.Range("O46").Value = "Normal spread capture: 198 bps"
$wsMM.Range("Q46").Value = 0

# ---------------------------------------------------------------------
# 5) EMAArbitrage sheet - close out the mirrored trade in row 8.
# ---------------------------------------------------------------------
$wsEMA = $wb.Worksheets.Item("EMAArbitrage")

$wsEMA.Range("G8").Value = 0.01
$wsEMA.Range("H8").Value = "CLOSED"
$wsEMA.Range("K8").Value = 100.27
$wsEMA.Range("P8").Value = "early_exit"
$wsEMA.Range("Q8").Value = 0.09
